$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Font size change -------------------------------------------------
# The title cell (C1, merged across C1:G2) uses a bold/italic Times New
# Roman font that was sized 25pt -> shrink it down to 20pt.
$ws.Range("C1").Font.Size = 20

# --- Column width changes ---------------------------------------------
# Previously only columns B and G had an explicit width (20 "characters"
# each). Now columns A, B, E, F and G all get explicit widths; B/G get
# narrower (17) and A/E/F become newly-narrow helper columns (6/5/6).
#
# Excel's ColumnWidth property is expressed in "characters" of the
# workbook's Normal-style font and includes a small built-in cell-padding
# offset before it is stored as the raw <col width="..."/> unit in the
# file, so we subtract that fixed padding (~0.8333 chars) back out here
# to land on the exact stored widths of 6 / 17 / 5 / 6 / 17.
$padding = 0.8333333333333334

$ws.Columns.Item(1).ColumnWidth = 6 - $padding
$ws.Columns.Item(2).ColumnWidth = 17 - $padding
$ws.Columns.Item(5).ColumnWidth = 5 - $padding
$ws.Columns.Item(6).ColumnWidth = 6 - $padding
$ws.Columns.Item(7).ColumnWidth = 17 - $padding
